# Auto-generated script to update market-price derived columns (H-N)
# across multiple worksheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(64, 8).Value = 24281.75  # H64: was 24323.459
$ws.Cells.Item(64, 9).Value = 3508.6667  # I64: was 3526.8823
$ws.Cells.Item(64, 10).Value = 86601  # J64: was 74829.42999999999
$ws.Cells.Item(64, 11).Value = 3508.6667  # K64: was 3526.8823
$ws.Cells.Item(64, 12).Value = 86601  # L64: was 74829.42999999999
$ws.Cells.Item(64, 13).Value = -3260.6667  # M64: was -3278.8823
$ws.Cells.Item(64, 14).Value = -87097  # N64: was -75325.42999999999
$ws.Cells.Item(67, 8).Value = 24281.75  # H67: was 24323.459
$ws.Cells.Item(67, 9).Value = 3508.6667  # I67: was 3526.8823
$ws.Cells.Item(67, 10).Value = 86601  # J67: was 74829.42999999999
$ws.Cells.Item(67, 11).Value = 3508.6667  # K67: was 3526.8823
$ws.Cells.Item(67, 12).Value = 86601  # L67: was 74829.42999999999
$ws.Cells.Item(67, 13).Value = -2650.6667  # M67: was -2668.8823
$ws.Cells.Item(67, 14).Value = -88317  # N67: was -76545.42999999999
$ws.Cells.Item(80, 8).Value = 4329769.5  # H80: was 5682797.5
$ws.Cells.Item(80, 9).Value = 251.93333  # I80: was 364.33334
$ws.Cells.Item(80, 10).Value = 15153564  # J80: was 22730098
$ws.Cells.Item(80, 11).Value = 755.79999  # K80: was 1093.00002
$ws.Cells.Item(80, 12).Value = 45460692  # L80: was 68190294
$ws.Cells.Item(80, 13).Value = 242.20001  # M80: was -95.00001999999995
$ws.Cells.Item(80, 14).Value = -45462688  # N80: was -68192290
$ws.Cells.Item(83, 8).Value = 4329769.5  # H83: was 5682797.5
$ws.Cells.Item(83, 9).Value = 251.93333  # I83: was 364.33334
$ws.Cells.Item(83, 10).Value = 15153564  # J83: was 22730098
$ws.Cells.Item(83, 11).Value = 2267.39997  # K83: was 3279.00006
$ws.Cells.Item(83, 12).Value = 136382076  # L83: was 204570882
$ws.Cells.Item(83, 13).Value = 2724.60003  # M83: was 1712.99994
$ws.Cells.Item(83, 14).Value = -136392060  # N83: was -204580866
# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(32, 8).Value = 5953.42  # H32: was 3726.52
$ws.Cells.Item(32, 9).Value = 5416.0947  # I32: was 3631
$ws.Cells.Item(32, 10).Value = 16162.6  # J32: was 8407
$ws.Cells.Item(32, 11).Value = 5416.0947  # K32: was 3631
$ws.Cells.Item(32, 12).Value = 16162.6  # L32: was 8407
$ws.Cells.Item(32, 13).Value = -5129.0947  # M32: was -3344
$ws.Cells.Item(32, 14).Value = -16736.6  # N32: was -8981
$ws.Cells.Item(42, 8).Value = 82531  # H42: was 0
$ws.Cells.Item(42, 10).Value = 82531  # J42: was 0
$ws.Cells.Item(42, 12).Value = 82531  # L42: was 0
$ws.Cells.Item(42, 14).Value = -83503  # N42: was None
$ws.Cells.Item(88, 8).Value = 11604.05  # H88: was 3109297
$ws.Cells.Item(88, 9).Value = 1373.25  # I88: was 1539.8572
$ws.Cells.Item(88, 10).Value = 18424.584  # J88: was 5284727
$ws.Cells.Item(88, 11).Value = 1373.25  # K88: was 1539.8572
$ws.Cells.Item(88, 12).Value = 18424.584  # L88: was 5284727
$ws.Cells.Item(88, 13).Value = -967.25  # M88: was -1133.8572
$ws.Cells.Item(88, 14).Value = -19236.584  # N88: was -5285539
$ws.Cells.Item(91, 8).Value = 11604.05  # H91: was 3109297
$ws.Cells.Item(91, 9).Value = 1373.25  # I91: was 1539.8572
$ws.Cells.Item(91, 10).Value = 18424.584  # J91: was 5284727
$ws.Cells.Item(91, 11).Value = 1373.25  # K91: was 1539.8572
$ws.Cells.Item(91, 12).Value = 18424.584  # L91: was 5284727
$ws.Cells.Item(91, 13).Value = 30.75  # M91: was -135.8571999999999
$ws.Cells.Item(91, 14).Value = -21232.584  # N91: was -5287535
$ws.Cells.Item(122, 8).Value = 1365.2858  # H122: was 1301.5
$ws.Cells.Item(122, 9).Value = 1271.4  # I122: was 1279.8334
$ws.Cells.Item(122, 10).Value = 1600  # J122: was 1317.75
$ws.Cells.Item(122, 11).Value = 3814.2  # K122: was 3839.5002
$ws.Cells.Item(122, 12).Value = 4800  # L122: was 3953.25
$ws.Cells.Item(122, 13).Value = -1364.2  # M122: was -1389.5002
$ws.Cells.Item(122, 14).Value = -9700  # N122: was -8853.25
# --- Sheet BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(10, 8).Value = 19335  # H10: was 25000
$ws.Cells.Item(10, 9).Value = 3005  # I10: was 0
$ws.Cells.Item(10, 10).Value = 27500  # J10: was 25000
$ws.Cells.Item(10, 11).Value = 3005  # K10: was 0
$ws.Cells.Item(10, 12).Value = 27500  # L10: was 25000
$ws.Cells.Item(10, 13).Value = -2865  # M10: was None
$ws.Cells.Item(10, 14).Value = -27780  # N10: was -25280
$ws.Cells.Item(134, 8).Value = 30335096  # H134: was 30335136
$ws.Cells.Item(134, 9).Value = 1372.2307  # I134: was 1551
$ws.Cells.Item(134, 10).Value = 143003220  # J134: was 125127590
$ws.Cells.Item(134, 11).Value = 4116.6921  # K134: was 4653
$ws.Cells.Item(134, 12).Value = 429009660  # L134: was 375382770
$ws.Cells.Item(134, 13).Value = -1581.6921  # M134: was -2118
$ws.Cells.Item(134, 14).Value = -429014730  # N134: was -375387840
# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)

$ws.Cells.Item(62, 8).Value = 2451  # H62: was 2500
$ws.Cells.Item(62, 9).Value = 2566.6667  # I62: was 2500
$ws.Cells.Item(62, 10).Value = 2335.3333  # J62: was 2500
$ws.Cells.Item(62, 11).Value = 2566.6667  # K62: was 2500
$ws.Cells.Item(62, 12).Value = 2335.3333  # L62: was 2500
$ws.Cells.Item(62, 13).Value = -1942.6667  # M62: was -1876
$ws.Cells.Item(62, 14).Value = -3583.3333  # N62: was -3748
$ws.Cells.Item(65, 8).Value = 2451  # H65: was 2500
$ws.Cells.Item(65, 9).Value = 2566.6667  # I65: was 2500
$ws.Cells.Item(65, 10).Value = 2335.3333  # J65: was 2500
$ws.Cells.Item(65, 11).Value = 12833.3335  # K65: was 12500
$ws.Cells.Item(65, 12).Value = 11676.6665  # L65: was 12500
$ws.Cells.Item(65, 13).Value = -9713.333500000001  # M65: was -9380
$ws.Cells.Item(65, 14).Value = -17916.6665  # N65: was -18740
$ws.Cells.Item(132, 8).Value = 35789.484  # H132: was 46957.227
$ws.Cells.Item(132, 9).Value = 50904.85  # I132: was 59914.65
$ws.Cells.Item(132, 10).Value = 2199.7778  # J132: was 2902
$ws.Cells.Item(132, 11).Value = 152714.55  # K132: was 179743.95
$ws.Cells.Item(132, 12).Value = 6599.3334  # L132: was 8706
$ws.Cells.Item(132, 13).Value = -150184.55  # M132: was -177213.95
$ws.Cells.Item(132, 14).Value = -11659.3334  # N132: was -13766
# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)

$ws.Cells.Item(59, 8).Value = 22735648  # H59: was 45469548
$ws.Cells.Item(59, 9).Value = 500  # I59: was 0
$ws.Cells.Item(59, 10).Value = 30314032  # J59: was 45469548
$ws.Cells.Item(59, 11).Value = 1500  # K59: was 0
$ws.Cells.Item(59, 12).Value = 90942096  # L59: was 136408644
$ws.Cells.Item(59, 13).Value = -960  # M59: was None
$ws.Cells.Item(59, 14).Value = -90943176  # N59: was -136409724
$ws.Cells.Item(113, 8).Value = 501.33334  # H113: was 402.76923
$ws.Cells.Item(113, 9).Value = 0  # I113: was 385.64285
$ws.Cells.Item(113, 10).Value = 501.33334  # J113: was 422.75
$ws.Cells.Item(113, 11).Value = 0  # K113: was 1156.92855
$ws.Cells.Item(113, 12).Value = 1504.00002  # L113: was 1268.25
$ws.Cells.Item(113, 13).ClearContents()  # M113: was 1013.07145
$ws.Cells.Item(113, 14).Value = -5844.000019999999  # N113: was -5608.25
# --- Sheet GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)

$ws.Cells.Item(29, 8).Value = 72004  # H29: was 74504
$ws.Cells.Item(29, 10).Value = 72004  # J29: was 74504
$ws.Cells.Item(29, 12).Value = 72004  # L29: was 74504
$ws.Cells.Item(29, 14).Value = -72584  # N29: was -75084
$ws.Cells.Item(102, 8).Value = 6973.143  # H102: was 9282.4
$ws.Cells.Item(102, 9).Value = 7802  # I102: was 11103
$ws.Cells.Item(102, 11).Value = 7802  # K102: was 11103
$ws.Cells.Item(102, 13).Value = -6180  # M102: was -9481
$ws.Cells.Item(122, 8).Value = 2636.182  # H122: was 2025.05
$ws.Cells.Item(122, 9).Value = 3271.1428  # I122: was 2224.5833
$ws.Cells.Item(122, 10).Value = 1525  # J122: was 1725.75
$ws.Cells.Item(122, 11).Value = 9813.428400000001  # K122: was 6673.749899999999
$ws.Cells.Item(122, 12).Value = 4575  # L122: was 5177.25
$ws.Cells.Item(122, 13).Value = -7363.428400000001  # M122: was -4223.749899999999
$ws.Cells.Item(122, 14).Value = -9475  # N122: was -10077.25
# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)

$ws.Cells.Item(7, 8).Value = 2402.7407  # H7: was 2304.6206
$ws.Cells.Item(7, 9).Value = 2594.389  # I7: was 2432.95
$ws.Cells.Item(7, 11).Value = 2594.389  # K7: was 2432.95
$ws.Cells.Item(7, 13).Value = -2482.389  # M7: was -2320.95
$ws.Cells.Item(26, 8).Value = 29800  # H26: was 9.5
$ws.Cells.Item(26, 9).Value = 0  # I26: was 9.5
$ws.Cells.Item(26, 10).Value = 29800  # J26: was 0
$ws.Cells.Item(26, 11).Value = 0  # K26: was 9.5
$ws.Cells.Item(26, 12).Value = 29800  # L26: was 0
$ws.Cells.Item(26, 13).ClearContents()  # M26: was 285.5
$ws.Cells.Item(26, 14).Value = -30390  # N26: was None
$ws.Cells.Item(40, 8).Value = 4632.143  # H40: was 61461.766
$ws.Cells.Item(40, 9).Value = 0  # I40: was 750
$ws.Cells.Item(40, 10).Value = 4632.143  # J40: was 69556.664
$ws.Cells.Item(40, 11).Value = 0  # K40: was 750
$ws.Cells.Item(40, 12).Value = 4632.143  # L40: was 69556.664
$ws.Cells.Item(40, 13).ClearContents()  # M40: was -614
$ws.Cells.Item(40, 14).Value = -4904.143  # N40: was -69828.664
$ws.Cells.Item(93, 8).Value = 1426.48  # H93: was 1397.7693
$ws.Cells.Item(93, 9).Value = 1052.7142  # I93: was 1027.8667
$ws.Cells.Item(93, 11).Value = 1052.7142  # K93: was 1027.8667
$ws.Cells.Item(93, 13).Value = 195.2858000000001  # M93: was 220.1333
$ws.Cells.Item(122, 8).Value = 2463.3333  # H122: was 2793.3333
$ws.Cells.Item(122, 9).Value = 2360  # I122: was 2766.6667
$ws.Cells.Item(122, 11).Value = 7080  # K122: was 8300.000100000001
$ws.Cells.Item(122, 13).Value = -4630  # M122: was -5850.000100000001
$ws.Cells.Item(126, 8).Value = 2402.7407  # H126: was 2304.6206
$ws.Cells.Item(126, 9).Value = 2594.389  # I126: was 2432.95
$ws.Cells.Item(126, 11).Value = 7783.167  # K126: was 7298.849999999999
$ws.Cells.Item(126, 13).Value = -5313.167  # M126: was -4828.849999999999
$ws.Cells.Item(132, 8).Value = 635601.1  # H132: was 309954.6
$ws.Cells.Item(132, 9).Value = 227252.44  # I132: was 70933.8
$ws.Cells.Item(132, 11).Value = 681757.3200000001  # K132: was 212801.4
$ws.Cells.Item(132, 13).Value = -679227.3200000001  # M132: was -210271.4
# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)

$ws.Cells.Item(62, 8).Value = 4888.7  # H62: was 4484.96
$ws.Cells.Item(62, 9).Value = 3464  # I62: was 3373
$ws.Cells.Item(62, 10).Value = 5140.1177  # J62: was 4696.7617
$ws.Cells.Item(62, 11).Value = 3464  # K62: was 3373
$ws.Cells.Item(62, 12).Value = 5140.1177  # L62: was 4696.7617
$ws.Cells.Item(62, 13).Value = -2840  # M62: was -2749
$ws.Cells.Item(62, 14).Value = -6388.1177  # N62: was -5944.7617
$ws.Cells.Item(65, 8).Value = 4888.7  # H65: was 4484.96
$ws.Cells.Item(65, 9).Value = 3464  # I65: was 3373
$ws.Cells.Item(65, 10).Value = 5140.1177  # J65: was 4696.7617
$ws.Cells.Item(65, 11).Value = 17320  # K65: was 16865
$ws.Cells.Item(65, 12).Value = 25700.5885  # L65: was 23483.8085
$ws.Cells.Item(65, 13).Value = -14200  # M65: was -13745
$ws.Cells.Item(65, 14).Value = -31940.5885  # N65: was -29723.8085
$ws.Cells.Item(81, 8).Value = 1582.95  # H81: was 1911.8148
$ws.Cells.Item(81, 9).Value = 1744.2222  # I81: was 2299.8333
$ws.Cells.Item(81, 10).Value = 1451  # J81: was 1601.4
$ws.Cells.Item(81, 11).Value = 3488.4444  # K81: was 4599.6666
$ws.Cells.Item(81, 12).Value = 2902  # L81: was 3202.8
$ws.Cells.Item(81, 13).Value = -2427.4444  # M81: was -3538.6666
$ws.Cells.Item(81, 14).Value = -5024  # N81: was -5324.8
$ws.Cells.Item(84, 8).Value = 1582.95  # H84: was 1911.8148
$ws.Cells.Item(84, 9).Value = 1744.2222  # I84: was 2299.8333
$ws.Cells.Item(84, 10).Value = 1451  # J84: was 1601.4
$ws.Cells.Item(84, 11).Value = 17442.222  # K84: was 22998.333
$ws.Cells.Item(84, 12).Value = 14510  # L84: was 16014
$ws.Cells.Item(84, 13).Value = -12138.222  # M84: was -17694.333
$ws.Cells.Item(84, 14).Value = -25118  # N84: was -26622
$ws.Cells.Item(132, 8).Value = 6933.6665  # H132: was 2556.0454
$ws.Cells.Item(132, 9).Value = 1391.1818  # I132: was 734.4286
$ws.Cells.Item(132, 10).Value = 13030.4  # J132: was 5743.875
$ws.Cells.Item(132, 11).Value = 4173.5454  # K132: was 2203.2858
$ws.Cells.Item(132, 12).Value = 39091.2  # L132: was 17231.625
$ws.Cells.Item(132, 13).Value = -1643.5454  # M132: was 326.7142000000003
$ws.Cells.Item(132, 14).Value = -44151.2  # N132: was -22291.625
$ws.Cells.Item(136, 8).Value = 1567962.8  # H136: was 1734754.2
$ws.Cells.Item(136, 9).Value = 2234716.8  # I136: was 2648530.8
$ws.Cells.Item(136, 11).Value = 6704150.399999999  # K136: was 7945592.399999999
$ws.Cells.Item(136, 13).Value = -6701600.399999999  # M136: was -7943042.399999999
